# Perturbation: "Changed sheets in perturbation_tests/to_be_reformatted/math_L_curve"
#
# Summary of the edit (derived from the OOXML diff):
#   - optimization_parameters sheet:
#       * header row 1: drop the redundant C1:F1 "value" labels (keep only A1/B1)
#       * row "Model" / "Sigmoid"  -> relabel A-cell to "production_function"
#       * a brand new row is inserted right after it: "L_curve" / 1
#       * the old "Deletion" row (under the "Sheet" row) is removed entirely
#       * this sheet becomes the active/selected sheet, with selection E33
#   - dcin5_log2_expression sheet: no longer the active/tab-selected sheet
#   - workbook: active tab moves from dcin5_log2_expression (index 3) to
#     optimization_parameters (index 6)

$wb = $excel.ActiveWorkbook

$params = $wb.Worksheets.Item("optimization_parameters")
$dcin5Expr = $wb.Worksheets.Item("dcin5_log2_expression")

# --- optimization_parameters: trim the duplicated header cells (C1:F1) ---
$params.Range("C1:F1").ClearContents()

# --- rename the "Model" row label ---
$params.Range("A8").Value = "production_function"

# --- insert the new L_curve row right below it ---
$params.Rows.Item(9).Insert()
$params.Range("A9").Value = "L_curve"
$params.Range("B9").Value = 1
$params.Range("B9").NumberFormat = "0.00E+00"

# --- remove the old "Deletion" row (now pushed down to row 17) ---
$params.Rows.Item(17).Delete()

# --- move the active sheet / selection to optimization_parameters ---
$dcin5Expr.Select()
$dcin5Expr.Range("B1:M1").Select()

$params.Select()
$params.Range("E33").Select()

$wb.Worksheets.Item("optimization_parameters").Activate()
